$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.386.10'
$ws.Range("D3").Value = '1.845.92'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.70'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6380'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07561'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2967'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.77'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07739'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.869.75'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.995'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6845'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.17'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009946'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.179'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("D18").Value = '29.412.60'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.61'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.48'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.579'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1408'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.394'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.67'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.468'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05711'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.01%  '
$ws.Range("E30").Value = '  -2.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.136'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.033'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.849'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.157'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7172'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.592'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = '1.254.53'
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9086'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.153'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9995'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.61'
$ws.Range("D43").ClearFormats()
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.065'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.82%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.162'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4032'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.705'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1128'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05741'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4628'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.16%  '
